$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.282.04"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.867.55"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'234.71"
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "'0.4700"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").Value = "'0.2855"
$ws.Range("E8").Value = "  -1.46%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  -2.37%  "

$ws.Range("D11").Value = "'0.07824"
$ws.Range("E11").Value = "  -1.59%  "

$ws.Range("D12").Value = "'96.80"
$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("D13").Value = "1.848.14"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").Value = "'0.6944"
$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").Value = "'5.084"
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").Value = "'268.69"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("D17").Value = "30.238.18"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").Value = "'0.000007692"
$ws.Range("E19").Value = "  +2.45%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "2.103.52"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").Value = "'1.001"

$ws.Range("D23").Value = "'5.251"
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("D24").Value = "'6.153"

$ws.Range("D25").Value = "'9.583"
$ws.Range("E25").Value = "  +4.26%  "

$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("E29").Value = "  -2.34%  "

$ws.Range("D30").Value = "'0.09892"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").Value = "'1.458"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").Value = "'4.049"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").Value = "'0.04728"
$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").Value = "'0.7028"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").Value = "'0.01875"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'2.768"
$ws.Range("E39").Value = "  +5.46%  "

$ws.Range("D40").Value = "'6.325"
$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("D41").Value = "'72.87"
$ws.Range("E41").Value = "  -1.40%  "

$ws.Range("D42").Value = "'1.949"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").Value = "'0.4167"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'0.8356"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "'103.13"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("D47").Value = "'977.21"
$ws.Range("E47").Value = "  +3.43%  "

$ws.Range("D49").Value = "'9.154"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "'34.53"

$ws.Range("D51").Value = "'0.05682"
$ws.Range("E51").Value = "  +0.30%  "
